# Fruta / hortaliza, semanal
# Insert a new weekly record at row 31 (shifting all subsequent rows down by
# one) on the "Pomelo" sheet, then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 31; this shifts rows 31..77 down
# to 32..78 (carrying all of their existing data/formatting with them).
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly observation.
$ws.Cells.Item(31, 1).Value = 9
$ws.Cells.Item(31, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(31, 3).Value = "Metropolitana"
$ws.Cells.Item(31, 4).Value = 44775
$ws.Cells.Item(31, 5).Value = 13
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100102
$ws.Cells.Item(31, 8).Value = "Cítricos"
$ws.Cells.Item(31, 9).Value = 100102006
$ws.Cells.Item(31, 10).Value = "Pomelo"
$ws.Cells.Item(31, 11).Value = "Start Ruby"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 280
$ws.Cells.Item(31, 14).Value = 8500
$ws.Cells.Item(31, 15).Value = 8500
$ws.Cells.Item(31, 16).Value = 8500
$ws.Cells.Item(31, 17).Value = "`$/caja 14 kilos"
$ws.Cells.Item(31, 18).Value = "Región Metropolitana"
$ws.Cells.Item(31, 19).Value = 607
$ws.Cells.Item(31, 20).Value = 14
